$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "tenant_id" column before the existing "update_usr_id" column (H),
# shifting update_usr_id -> I and update_time -> J.
$ws.Range("H1").EntireColumn.Insert()

# New header comment-row cell (row 1) describing the tenant_id field + its
# dropdown data-validation generation, mirroring the other *_lbl header cells.
$ws.Range("H1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

# New sample-row cell (row 2) rendering the tenant_id label for each model row.
$ws.Range("H2").Value = '<%=model.tenant_id_lbl%>'
